$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# Column E ("authority") values change from "ks.gov" to "kshs.org" for all data rows.
$rng = $ws.Range("E2:E28")
$rng.Value = "kshs.org"

# Update the window scroll position.
$wb.Windows.Item(1).ScrollRow = 1

# Update selection to E2:E28 with active cell E2.
$ws.Range("E2:E28").Select()
